# Auto-generated Excel COM-interop script applying the Lamia_Profits data refresh diff.
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H:N) across all 8 crafting-job sheets
# to reflect the latest scheduled market-data pull.
$wb = $excel.ActiveWorkbook

# ===== ALC =====
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(15, 8).Value = 1428.5223
$ws.Cells.Item(15, 9).Value = 1428.5223
$ws.Cells.Item(15, 11).Value = 4285.5669
$ws.Cells.Item(15, 13).Value = -4116.5669
$ws.Cells.Item(19, 8).Value = 721
$ws.Cells.Item(19, 9).Value = 537.5
$ws.Cells.Item(19, 10).Value = 1149.1666
$ws.Cells.Item(19, 11).Value = 537.5
$ws.Cells.Item(19, 12).Value = 1149.1666
$ws.Cells.Item(19, 13).Value = -362.5
$ws.Cells.Item(19, 14).Value = -1499.1666
$ws.Cells.Item(40, 8).Value = 4912.364
$ws.Cells.Item(40, 9).Value = 1715
$ws.Cells.Item(40, 10).Value = 8749.200000000001
$ws.Cells.Item(40, 11).Value = 1715
$ws.Cells.Item(40, 12).Value = 8749.200000000001
$ws.Cells.Item(40, 13).Value = -1540
$ws.Cells.Item(40, 14).Value = -9099.200000000001
$ws.Cells.Item(62, 8).Value = 8450.611000000001
$ws.Cells.Item(62, 9).Value = 5524.75
$ws.Cells.Item(62, 10).Value = 9286.571
$ws.Cells.Item(62, 11).Value = 5524.75
$ws.Cells.Item(62, 12).Value = 9286.571
$ws.Cells.Item(62, 13).Value = -4900.75
$ws.Cells.Item(62, 14).Value = -10534.571
$ws.Cells.Item(65, 8).Value = 8450.611000000001
$ws.Cells.Item(65, 9).Value = 5524.75
$ws.Cells.Item(65, 10).Value = 9286.571
$ws.Cells.Item(65, 11).Value = 27623.75
$ws.Cells.Item(65, 12).Value = 46432.855
$ws.Cells.Item(65, 13).Value = -24503.75
$ws.Cells.Item(65, 14).Value = -52672.855
$ws.Cells.Item(74, 8).Value = 18001
$ws.Cells.Item(74, 9).Value = 17666.666
$ws.Cells.Item(74, 10).Value = 19004
$ws.Cells.Item(74, 11).Value = 17666.666
$ws.Cells.Item(74, 12).Value = 19004
$ws.Cells.Item(74, 13).Value = -16730.666
$ws.Cells.Item(74, 14).Value = -20876
$ws.Cells.Item(76, 8).Value = 8578.380999999999
$ws.Cells.Item(76, 9).Value = 8013.8
$ws.Cells.Item(76, 10).Value = 9091.637000000001
$ws.Cells.Item(76, 11).Value = 8013.8
$ws.Cells.Item(76, 12).Value = 9091.637000000001
$ws.Cells.Item(76, 13).Value = -7698.8
$ws.Cells.Item(76, 14).Value = -9721.637000000001
$ws.Cells.Item(77, 8).Value = 18001
$ws.Cells.Item(77, 9).Value = 17666.666
$ws.Cells.Item(77, 10).Value = 19004
$ws.Cells.Item(77, 11).Value = 88333.33
$ws.Cells.Item(77, 12).Value = 95020
$ws.Cells.Item(77, 13).Value = -83653.33
$ws.Cells.Item(77, 14).Value = -104380
$ws.Cells.Item(79, 8).Value = 8578.380999999999
$ws.Cells.Item(79, 9).Value = 8013.8
$ws.Cells.Item(79, 10).Value = 9091.637000000001
$ws.Cells.Item(79, 11).Value = 8013.8
$ws.Cells.Item(79, 12).Value = 9091.637000000001
$ws.Cells.Item(79, 13).Value = -6921.8
$ws.Cells.Item(79, 14).Value = -11275.637
$ws.Cells.Item(86, 8).Value = 4403.5884
$ws.Cells.Item(86, 9).Value = 4116.6665
$ws.Cells.Item(86, 10).Value = 4560.091
$ws.Cells.Item(86, 11).Value = 4116.6665
$ws.Cells.Item(86, 12).Value = 4560.091
$ws.Cells.Item(86, 13).Value = -2993.6665
$ws.Cells.Item(86, 14).Value = -6806.091
$ws.Cells.Item(89, 8).Value = 4403.5884
$ws.Cells.Item(89, 9).Value = 4116.6665
$ws.Cells.Item(89, 10).Value = 4560.091
$ws.Cells.Item(89, 11).Value = 20583.3325
$ws.Cells.Item(89, 12).Value = 22800.455
$ws.Cells.Item(89, 13).Value = -14967.3325
$ws.Cells.Item(89, 14).Value = -34032.455
$ws.Cells.Item(92, 8).Value = 2553.56
$ws.Cells.Item(92, 9).Value = 1423.5238
$ws.Cells.Item(92, 11).Value = 1423.5238
$ws.Cells.Item(92, 13).Value = -175.5237999999999
$ws.Cells.Item(98, 8).Value = 274615.38
$ws.Cells.Item(98, 9).Value = 889.05884
$ws.Cells.Item(98, 10).Value = 1205284.8
$ws.Cells.Item(98, 11).Value = 889.05884
$ws.Cells.Item(98, 12).Value = 1205284.8
$ws.Cells.Item(98, 13).Value = 608.94116
$ws.Cells.Item(98, 14).Value = -1208280.8
$ws.Cells.Item(106, 8).Value = 5256.615
$ws.Cells.Item(106, 9).Value = 5226.5
$ws.Cells.Item(106, 10).Value = 5282.4287
$ws.Cells.Item(106, 11).Value = 5226.5
$ws.Cells.Item(106, 12).Value = 5282.4287
$ws.Cells.Item(106, 13).Value = -4595.5
$ws.Cells.Item(106, 14).Value = -6544.4287
$ws.Cells.Item(122, 8).Value = 274615.38
$ws.Cells.Item(122, 9).Value = 889.05884
$ws.Cells.Item(122, 10).Value = 1205284.8
$ws.Cells.Item(122, 11).Value = 2667.17652
$ws.Cells.Item(122, 12).Value = 3615854.4
$ws.Cells.Item(122, 13).Value = -217.17652
$ws.Cells.Item(122, 14).Value = -3620754.4
$ws.Cells.Item(132, 8).Value = 1475.2407
$ws.Cells.Item(132, 9).Value = 916.0417
$ws.Cells.Item(132, 11).Value = 2748.1251
$ws.Cells.Item(132, 13).Value = -218.1251000000002
$ws.Cells.Item(135, 8).Value = 2345.9443
$ws.Cells.Item(135, 9).Value = 1895.7059
$ws.Cells.Item(135, 11).Value = 17061.3531
$ws.Cells.Item(135, 13).Value = -14526.3531
$ws.Cells.Item(137, 8).Value = 11630655
$ws.Cells.Item(137, 9).Value = 35716150
$ws.Cells.Item(137, 10).Value = 3173.5344
$ws.Cells.Item(137, 11).Value = 107148450
$ws.Cells.Item(137, 12).Value = 9520.6032
$ws.Cells.Item(137, 13).Value = -107145900
$ws.Cells.Item(137, 14).Value = -14620.6032
$ws.Cells.Item(138, 8).Value = 3844.4814
$ws.Cells.Item(138, 10).Value = 4151.778
$ws.Cells.Item(138, 12).Value = 12455.334
$ws.Cells.Item(138, 14).Value = -22735.334
$ws.Cells.Item(139, 8).Value = 48171.547
$ws.Cells.Item(139, 10).Value = 48171.547
$ws.Cells.Item(139, 12).Value = 48171.547
$ws.Cells.Item(139, 14).Value = -58451.547
$ws.Cells.Item(141, 8).Value = 2973.818
$ws.Cells.Item(141, 9).Value = 2448.4375
$ws.Cells.Item(141, 10).Value = 4374.8335
$ws.Cells.Item(141, 11).Value = 7345.3125
$ws.Cells.Item(141, 12).Value = 13124.5005
$ws.Cells.Item(141, 13).Value = -2165.3125
$ws.Cells.Item(141, 14).Value = -23484.5005

# ===== ARM =====
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(5, 8).Value = 288.58334
$ws.Cells.Item(5, 9).Value = 173.6
$ws.Cells.Item(5, 11).Value = 173.6
$ws.Cells.Item(5, 13).Value = -61.59999999999999
$ws.Cells.Item(32, 8).Value = 4744.4346
$ws.Cells.Item(32, 9).Value = 4247.121
$ws.Cells.Item(32, 10).Value = 50000
$ws.Cells.Item(32, 11).Value = 4247.121
$ws.Cells.Item(32, 12).Value = 50000
$ws.Cells.Item(32, 13).Value = -3960.121
$ws.Cells.Item(32, 14).Value = -50574
$ws.Cells.Item(74, 8).Value = 2834.625
$ws.Cells.Item(74, 9).Value = 2439.8696
$ws.Cells.Item(74, 11).Value = 2439.8696
$ws.Cells.Item(74, 13).Value = -1565.8696
$ws.Cells.Item(77, 8).Value = 2834.625
$ws.Cells.Item(77, 9).Value = 2439.8696
$ws.Cells.Item(77, 11).Value = 12199.348
$ws.Cells.Item(77, 13).Value = -7831.348
$ws.Cells.Item(88, 8).Value = 5774.25
$ws.Cells.Item(88, 9).Value = 6399.25
$ws.Cells.Item(88, 10).Value = 5565.9165
$ws.Cells.Item(88, 11).Value = 6399.25
$ws.Cells.Item(88, 12).Value = 5565.9165
$ws.Cells.Item(88, 13).Value = -5993.25
$ws.Cells.Item(88, 14).Value = -6377.9165
$ws.Cells.Item(91, 8).Value = 5774.25
$ws.Cells.Item(91, 9).Value = 6399.25
$ws.Cells.Item(91, 10).Value = 5565.9165
$ws.Cells.Item(91, 11).Value = 6399.25
$ws.Cells.Item(91, 12).Value = 5565.9165
$ws.Cells.Item(91, 13).Value = -4995.25
$ws.Cells.Item(91, 14).Value = -8373.916499999999
$ws.Cells.Item(110, 8).Value = 3343.625
$ws.Cells.Item(110, 10).Value = 10506.5
$ws.Cells.Item(110, 12).Value = 10506.5
$ws.Cells.Item(110, 14).Value = -14596.5
$ws.Cells.Item(122, 8).Value = 2878.973
$ws.Cells.Item(122, 9).Value = 2145.5
$ws.Cells.Item(122, 10).Value = 4233.077
$ws.Cells.Item(122, 11).Value = 6436.5
$ws.Cells.Item(122, 12).Value = 12699.231
$ws.Cells.Item(122, 13).Value = -3986.5
$ws.Cells.Item(122, 14).Value = -17599.231
$ws.Cells.Item(132, 8).Value = 2607
$ws.Cells.Item(132, 9).Value = 1786.65
$ws.Cells.Item(132, 11).Value = 5359.950000000001
$ws.Cells.Item(132, 13).Value = -2829.950000000001

# ===== BSM =====
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(4, 8).Value = 288.58334
$ws.Cells.Item(4, 9).Value = 173.6
$ws.Cells.Item(4, 11).Value = 173.6
$ws.Cells.Item(4, 13).Value = -58.59999999999999
$ws.Cells.Item(22, 8).Value = 225
$ws.Cells.Item(22, 9).Value = 250
$ws.Cells.Item(22, 10).Value = 175
$ws.Cells.Item(22, 11).Value = 250
$ws.Cells.Item(22, 12).Value = 175
$ws.Cells.Item(22, 13).Value = -77
$ws.Cells.Item(22, 14).Value = -521
$ws.Cells.Item(107, 8).Value = 626.1429000000001
$ws.Cells.Item(107, 9).Value = 564
$ws.Cells.Item(107, 11).Value = 564
$ws.Cells.Item(107, 13).Value = 1356
$ws.Cells.Item(134, 8).Value = 1712.08
$ws.Cells.Item(134, 9).Value = 1440.6123
$ws.Cells.Item(134, 11).Value = 4321.8369
$ws.Cells.Item(134, 13).Value = -1786.8369

# ===== CRP =====
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 33254.973
$ws.Cells.Item(31, 9).Value = 1348.9445
$ws.Cells.Item(31, 10).Value = 67037.82000000001
$ws.Cells.Item(31, 11).Value = 1348.9445
$ws.Cells.Item(31, 12).Value = 67037.82000000001
$ws.Cells.Item(31, 13).Value = -1053.9445
$ws.Cells.Item(31, 14).Value = -67627.82000000001
$ws.Cells.Item(34, 8).Value = 33254.973
$ws.Cells.Item(34, 9).Value = 1348.9445
$ws.Cells.Item(34, 10).Value = 67037.82000000001
$ws.Cells.Item(34, 11).Value = 1348.9445
$ws.Cells.Item(34, 12).Value = 67037.82000000001
$ws.Cells.Item(34, 13).Value = -1146.9445
$ws.Cells.Item(34, 14).Value = -67441.82000000001
$ws.Cells.Item(58, 8).Value = 2973.8462
$ws.Cells.Item(58, 9).Value = 1584.7858
$ws.Cells.Item(58, 10).Value = 6509.636
$ws.Cells.Item(58, 11).Value = 1584.7858
$ws.Cells.Item(58, 12).Value = 6509.636
$ws.Cells.Item(58, 13).Value = -1381.7858
$ws.Cells.Item(58, 14).Value = -6915.636
$ws.Cells.Item(99, 8).Value = 2636.5715
$ws.Cells.Item(99, 9).Value = 2545.182
$ws.Cells.Item(99, 10).Value = 2695.7058
$ws.Cells.Item(99, 11).Value = 2545.182
$ws.Cells.Item(99, 12).Value = 2695.7058
$ws.Cells.Item(99, 13).Value = -1047.182
$ws.Cells.Item(99, 14).Value = -5691.7058
$ws.Cells.Item(107, 8).Value = 1868.9048
$ws.Cells.Item(107, 9).Value = 1204.6666
$ws.Cells.Item(107, 11).Value = 1204.6666
$ws.Cells.Item(107, 13).Value = 715.3334
$ws.Cells.Item(120, 8).Value = 49999
$ws.Cells.Item(120, 10).Value = 49999
$ws.Cells.Item(120, 12).Value = 49999
$ws.Cells.Item(120, 14).Value = -57257
$ws.Cells.Item(122, 8).Value = 2869.9211
$ws.Cells.Item(122, 9).Value = 1225
$ws.Cells.Item(122, 10).Value = 5131.6875
$ws.Cells.Item(122, 11).Value = 3675
$ws.Cells.Item(122, 12).Value = 15395.0625
$ws.Cells.Item(122, 13).Value = -1225
$ws.Cells.Item(122, 14).Value = -20295.0625
$ws.Cells.Item(126, 8).Value = 2636.5715
$ws.Cells.Item(126, 9).Value = 2545.182
$ws.Cells.Item(126, 10).Value = 2695.7058
$ws.Cells.Item(126, 11).Value = 7635.545999999999
$ws.Cells.Item(126, 12).Value = 8087.117400000001
$ws.Cells.Item(126, 13).Value = -5165.545999999999
$ws.Cells.Item(126, 14).Value = -13027.1174
$ws.Cells.Item(132, 8).Value = 3107.125
$ws.Cells.Item(132, 9).Value = 2548.9312
$ws.Cells.Item(132, 10).Value = 8503
$ws.Cells.Item(132, 11).Value = 7646.7936
$ws.Cells.Item(132, 12).Value = 25509
$ws.Cells.Item(132, 13).Value = -5116.7936
$ws.Cells.Item(132, 14).Value = -30569
$ws.Cells.Item(134, 8).Value = 1676.2858
$ws.Cells.Item(134, 9).Value = 1350.9756
$ws.Cells.Item(134, 11).Value = 4052.9268
$ws.Cells.Item(134, 13).Value = -1517.9268
$ws.Cells.Item(136, 8).Value = 2973.8462
$ws.Cells.Item(136, 9).Value = 1584.7858
$ws.Cells.Item(136, 10).Value = 6509.636
$ws.Cells.Item(136, 11).Value = 4754.357400000001
$ws.Cells.Item(136, 12).Value = 19528.908
$ws.Cells.Item(136, 13).Value = -2204.357400000001
$ws.Cells.Item(136, 14).Value = -24628.908

# ===== CUL =====
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 30815434
$ws.Cells.Item(4, 9).Value = 67532584
$ws.Cells.Item(4, 10).Value = 7697229.5
$ws.Cells.Item(4, 11).Value = 202597752
$ws.Cells.Item(4, 12).Value = 23091688.5
$ws.Cells.Item(4, 13).Value = -202597640
$ws.Cells.Item(4, 14).Value = -23091912.5
$ws.Cells.Item(23, 8).Value = 407.77777
$ws.Cells.Item(23, 9).Value = 253.16667
$ws.Cells.Item(23, 10).Value = 485.08334
$ws.Cells.Item(23, 11).Value = 759.50001
$ws.Cells.Item(23, 12).Value = 1455.25002
$ws.Cells.Item(23, 13).Value = -524.50001
$ws.Cells.Item(23, 14).Value = -1925.25002
$ws.Cells.Item(81, 8).Value = 4387.8335
$ws.Cells.Item(81, 10).Value = 5976.6665
$ws.Cells.Item(81, 12).Value = 17929.9995
$ws.Cells.Item(81, 14).Value = -20175.9995
$ws.Cells.Item(84, 8).Value = 4387.8335
$ws.Cells.Item(84, 10).Value = 5976.6665
$ws.Cells.Item(84, 12).Value = 53789.9985
$ws.Cells.Item(84, 14).Value = -65021.9985
$ws.Cells.Item(107, 8).Value = 459906.97
$ws.Cells.Item(107, 9).Value = 283.26828
$ws.Cells.Item(107, 10).Value = 1157854.1
$ws.Cells.Item(107, 11).Value = 849.80484
$ws.Cells.Item(107, 12).Value = 3473562.3
$ws.Cells.Item(107, 13).Value = 1070.19516
$ws.Cells.Item(107, 14).Value = -3477402.3
$ws.Cells.Item(116, 8).Value = 5342100.5
$ws.Cells.Item(116, 9).Value = 8008135
$ws.Cells.Item(116, 10).Value = 4009083
$ws.Cells.Item(116, 11).Value = 24024405
$ws.Cells.Item(116, 12).Value = 12027249
$ws.Cells.Item(116, 13).Value = -24020963
$ws.Cells.Item(116, 14).Value = -12034133
$ws.Cells.Item(131, 8).Value = 4841842.5
$ws.Cells.Item(131, 10).Value = 3413211.8
$ws.Cells.Item(131, 12).Value = 10239635.4
$ws.Cells.Item(131, 14).Value = -10249715.4

# ===== GSM =====
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 561879.75
$ws.Cells.Item(80, 9).Value = 1668966.6
$ws.Cells.Item(80, 10).Value = 8336.333000000001
$ws.Cells.Item(80, 11).Value = 1668966.6
$ws.Cells.Item(80, 12).Value = 8336.333000000001
$ws.Cells.Item(80, 13).Value = -1667968.6
$ws.Cells.Item(80, 14).Value = -10332.333
$ws.Cells.Item(83, 8).Value = 561879.75
$ws.Cells.Item(83, 9).Value = 1668966.6
$ws.Cells.Item(83, 10).Value = 8336.333000000001
$ws.Cells.Item(83, 11).Value = 8344833
$ws.Cells.Item(83, 12).Value = 41681.665
$ws.Cells.Item(83, 13).Value = -8339841
$ws.Cells.Item(83, 14).Value = -51665.665
$ws.Cells.Item(86, 8).Value = 0
$ws.Cells.Item(86, 10).Value = 0
$ws.Cells.Item(86, 12).Value = 0
$ws.Cells.Item(86, 14).Value = ""
$ws.Cells.Item(89, 8).Value = 0
$ws.Cells.Item(89, 10).Value = 0
$ws.Cells.Item(89, 12).Value = 0
$ws.Cells.Item(89, 14).Value = ""
$ws.Cells.Item(97, 8).Value = 1591.8334
$ws.Cells.Item(97, 10).Value = 4218.8
$ws.Cells.Item(97, 12).Value = 4218.8
$ws.Cells.Item(97, 14).Value = -5210.8
$ws.Cells.Item(102, 8).Value = 3230.0454
$ws.Cells.Item(102, 9).Value = 1690.5
$ws.Cells.Item(102, 11).Value = 1690.5
$ws.Cells.Item(102, 13).Value = -68.5
$ws.Cells.Item(113, 8).Value = 3995.2
$ws.Cells.Item(113, 9).Value = 3367.25
$ws.Cells.Item(113, 11).Value = 3367.25
$ws.Cells.Item(113, 13).Value = -1197.25
$ws.Cells.Item(131, 8).Value = 47900
$ws.Cells.Item(131, 10).Value = 47900
$ws.Cells.Item(131, 12).Value = 47900
$ws.Cells.Item(131, 14).Value = -57980
$ws.Cells.Item(132, 8).Value = 2404.9846
$ws.Cells.Item(132, 9).Value = 2085.6453
$ws.Cells.Item(132, 10).Value = 9004.666999999999
$ws.Cells.Item(132, 11).Value = 6256.9359
$ws.Cells.Item(132, 12).Value = 27014.001
$ws.Cells.Item(132, 13).Value = -3726.9359
$ws.Cells.Item(132, 14).Value = -32074.001

# ===== LTW =====
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 3696.9524
$ws.Cells.Item(22, 9).Value = 1177.9166
$ws.Cells.Item(22, 10).Value = 7055.6665
$ws.Cells.Item(22, 11).Value = 1177.9166
$ws.Cells.Item(22, 12).Value = 7055.6665
$ws.Cells.Item(22, 13).Value = -882.9166
$ws.Cells.Item(22, 14).Value = -7645.6665
$ws.Cells.Item(27, 8).Value = 3696.9524
$ws.Cells.Item(27, 9).Value = 1177.9166
$ws.Cells.Item(27, 10).Value = 7055.6665
$ws.Cells.Item(27, 11).Value = 1177.9166
$ws.Cells.Item(27, 12).Value = 7055.6665
$ws.Cells.Item(27, 13).Value = -1070.9166
$ws.Cells.Item(27, 14).Value = -7269.6665
$ws.Cells.Item(40, 8).Value = 13895.643
$ws.Cells.Item(40, 9).Value = 32774.5
$ws.Cells.Item(40, 11).Value = 32774.5
$ws.Cells.Item(40, 13).Value = -32638.5
$ws.Cells.Item(43, 8).Value = 28911
$ws.Cells.Item(43, 9).Value = 28911
$ws.Cells.Item(43, 11).Value = 28911
$ws.Cells.Item(43, 13).Value = -28718
$ws.Cells.Item(61, 8).Value = 2903.56
$ws.Cells.Item(61, 9).Value = 2360.2273
$ws.Cells.Item(61, 10).Value = 6888
$ws.Cells.Item(61, 11).Value = 2360.2273
$ws.Cells.Item(61, 12).Value = 6888
$ws.Cells.Item(61, 13).Value = -2158.2273
$ws.Cells.Item(61, 14).Value = -7292
$ws.Cells.Item(68, 8).Value = 7425.7144
$ws.Cells.Item(68, 9).Value = 5349.5
$ws.Cells.Item(68, 10).Value = 7914.2354
$ws.Cells.Item(68, 11).Value = 5349.5
$ws.Cells.Item(68, 12).Value = 7914.2354
$ws.Cells.Item(68, 13).Value = -4600.5
$ws.Cells.Item(68, 14).Value = -9412.2354
$ws.Cells.Item(71, 8).Value = 7425.7144
$ws.Cells.Item(71, 9).Value = 5349.5
$ws.Cells.Item(71, 10).Value = 7914.2354
$ws.Cells.Item(71, 11).Value = 26747.5
$ws.Cells.Item(71, 12).Value = 39571.177
$ws.Cells.Item(71, 13).Value = -23003.5
$ws.Cells.Item(71, 14).Value = -47059.177
$ws.Cells.Item(93, 8).Value = 2277.6296
$ws.Cells.Item(93, 9).Value = 2260.087
$ws.Cells.Item(93, 11).Value = 2260.087
$ws.Cells.Item(93, 13).Value = -1012.087
$ws.Cells.Item(113, 8).Value = 2903.56
$ws.Cells.Item(113, 9).Value = 2360.2273
$ws.Cells.Item(113, 10).Value = 6888
$ws.Cells.Item(113, 11).Value = 2360.2273
$ws.Cells.Item(113, 12).Value = 6888
$ws.Cells.Item(113, 13).Value = -190.2273
$ws.Cells.Item(113, 14).Value = -11228
$ws.Cells.Item(122, 8).Value = 241303.64
$ws.Cells.Item(122, 9).Value = 312981.3
$ws.Cells.Item(122, 11).Value = 938943.8999999999
$ws.Cells.Item(122, 13).Value = -936493.8999999999
$ws.Cells.Item(132, 8).Value = 7297.375
$ws.Cells.Item(132, 9).Value = 5624.8335
$ws.Cells.Item(132, 10).Value = 12315
$ws.Cells.Item(132, 11).Value = 16874.5005
$ws.Cells.Item(132, 12).Value = 36945
$ws.Cells.Item(132, 13).Value = -14344.5005
$ws.Cells.Item(132, 14).Value = -42005

# ===== WVR =====
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(18, 8).Value = 28338
$ws.Cells.Item(18, 10).Value = 28007
$ws.Cells.Item(18, 12).Value = 28007
$ws.Cells.Item(18, 14).Value = -28353
$ws.Cells.Item(25, 8).Value = 0
$ws.Cells.Item(25, 10).Value = 0
$ws.Cells.Item(25, 12).Value = 0
$ws.Cells.Item(25, 14).Value = ""
$ws.Cells.Item(28, 8).Value = 24999.5
$ws.Cells.Item(28, 10).Value = 24999.5
$ws.Cells.Item(28, 12).Value = 24999.5
$ws.Cells.Item(28, 14).Value = -25695.5
$ws.Cells.Item(62, 8).Value = 8221.111000000001
$ws.Cells.Item(62, 9).Value = 7990
$ws.Cells.Item(62, 10).Value = 8250
$ws.Cells.Item(62, 11).Value = 7990
$ws.Cells.Item(62, 12).Value = 8250
$ws.Cells.Item(62, 13).Value = -7366
$ws.Cells.Item(62, 14).Value = -9498
$ws.Cells.Item(65, 8).Value = 8221.111000000001
$ws.Cells.Item(65, 9).Value = 7990
$ws.Cells.Item(65, 10).Value = 8250
$ws.Cells.Item(65, 11).Value = 39950
$ws.Cells.Item(65, 12).Value = 41250
$ws.Cells.Item(65, 13).Value = -36830
$ws.Cells.Item(65, 14).Value = -47490
$ws.Cells.Item(122, 8).Value = 2673.5483
$ws.Cells.Item(122, 9).Value = 1982.4584
$ws.Cells.Item(122, 11).Value = 5947.3752
$ws.Cells.Item(122, 13).Value = -3497.3752
$ws.Cells.Item(126, 8).Value = 3118.375
$ws.Cells.Item(126, 9).Value = 2608.2942
$ws.Cells.Item(126, 10).Value = 4357.143
$ws.Cells.Item(126, 11).Value = 7824.882599999999
$ws.Cells.Item(126, 12).Value = 13071.429
$ws.Cells.Item(126, 13).Value = -5354.882599999999
$ws.Cells.Item(126, 14).Value = -18011.429
$ws.Cells.Item(132, 8).Value = 2105.2
$ws.Cells.Item(132, 9).Value = 1522.4482
$ws.Cells.Item(132, 11).Value = 4567.3446
$ws.Cells.Item(132, 13).Value = -2037.3446
$ws.Cells.Item(136, 8).Value = 3114.7307
$ws.Cells.Item(136, 9).Value = 1304.2
$ws.Cells.Item(136, 10).Value = 9149.833000000001
$ws.Cells.Item(136, 11).Value = 3912.6
$ws.Cells.Item(136, 12).Value = 27449.499
$ws.Cells.Item(136, 13).Value = -1362.6
$ws.Cells.Item(136, 14).Value = -32549.499

Write-Host "Lamia_Profits data refresh applied."